$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "Li4SiO4-Be pebble bed" breeder rows: the original rows
# (15-22) all shared that one label. They now split into two distinct
# materials:
#   - rows 19-22 (the new full BISO/heterogeneous unit-test rows) keep
#     the original shared-string slot, renamed "Li4SiO4-Be-He"
#   - rows 15-18 get a brand-new label "Li4SiO4-Be"
# Doing the 19-22 rename first lets the engine reuse/rename the existing
# shared-string slot for "Li4SiO4-Be-He" and allocate a fresh slot for
# "Li4SiO4-Be" used by rows 15-18.
$ws.Range("A19").Value = "Li4SiO4-Be-He"
$ws.Range("A20").Value = "Li4SiO4-Be-He"
$ws.Range("A21").Value = "Li4SiO4-Be-He"
$ws.Range("A22").Value = "Li4SiO4-Be-He"

$ws.Range("A15").Value = "Li4SiO4-Be"
$ws.Range("A16").Value = "Li4SiO4-Be"
$ws.Range("A17").Value = "Li4SiO4-Be"
$ws.Range("A18").Value = "Li4SiO4-Be"

# Updated measured flux/reaction-rate inputs for the new full
# heterogeneous HCPB BISO run; dependent formulas (I/J/K columns)
# recalculate automatically.
$ws.Range("G19").Value = 2743.65
$ws.Range("G21").Value = 93396
$ws.Range("H21").Value = 101228

# Rows 20 and 22 are "Heterog"-only rows with no homogeneous (I) value,
# so their shared %-change formula used to divide by zero. Clear those
# results back to blank cells.
$ws.Range("K20").ClearContents()
$ws.Range("K22").ClearContents()

# Restore the saved cursor position/selection.
$ws.Range("G27").Select()
